$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Delete the row containing account 005092207 / BRUNO / 21486.75 (Excel row 3)
$ws.Rows.Item(3).Delete()
